# Retested some photos to make sure error handling works.
# The "Traceback ... UnboundLocalError" failure for the two photos that used to
# error out with that Python traceback now instead raises a clean
# "SystemExit: No faces were found" message, and the retest date moved from
# 4/25/2013 to 5/6/2013 (serial 41400). Row heights are adjusted to match the
# new (shorter) text, and the view is left positioned near the edited rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 - "https://onlinehealthsafe.com/..." entry
$ws.Range("C3").Value = "SystemExit: No faces were found"
$ws.Range("D3").Value = 41400
$ws.Rows.Item(3).RowHeight = 90

# Row 16 - "http://bloximages..." entry
$ws.Range("C16").Value = "SystemExit: No faces were found"
$ws.Range("D16").Value = 41400
$ws.Rows.Item(16).RowHeight = 120

# Update the view to reflect where work was last happening in the sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 14
$ws.Range("C17").Select()
